$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 13515958
$ws.Range("J112").Value = 14288002
$ws.Range("L112").Value = 42864006
$ws.Range("N112").Value = -42866222
$ws.Range("H135").Value = 449.36
$ws.Range("I135").Value = 418.86957
$ws.Range("K135").Value = 3769.82613
$ws.Range("M135").Value = -1234.82613
$ws.Range("H137").Value = 19356.217
$ws.Range("I137").Value = 19998.334
$ws.Range("K137").Value = 59995.00199999999
$ws.Range("M137").Value = -57445.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4288.397
$ws.Range("I32").Value = 4082.2458
$ws.Range("K32").Value = 4082.2458
$ws.Range("M32").Value = -3795.2458
$ws.Range("H61").Value = 4456.615
$ws.Range("I61").Value = 3812.4443
$ws.Range("J61").Value = 5906
$ws.Range("K61").Value = 3812.4443
$ws.Range("L61").Value = 5906
$ws.Range("M61").Value = -3600.4443
$ws.Range("N61").Value = -6330
$ws.Range("H126").Value = 6565
$ws.Range("I126").Value = 6565
$ws.Range("K126").Value = 19695
$ws.Range("M126").Value = -17225
$ws.Range("H132").Value = 29679.53
$ws.Range("I132").Value = 2079.1738
$ws.Range("K132").Value = 6237.5214
$ws.Range("M132").Value = -3707.5214
$ws.Range("H136").Value = 4456.615
$ws.Range("I136").Value = 3812.4443
$ws.Range("J136").Value = 5906
$ws.Range("K136").Value = 11437.3329
$ws.Range("L136").Value = 17718
$ws.Range("M136").Value = -8887.332900000001
$ws.Range("N136").Value = -22818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9870.559999999999
$ws.Range("I86").Value = 1638.0526
$ws.Range("K86").Value = 1638.0526
$ws.Range("M86").Value = -515.0526
$ws.Range("H89").Value = 9870.559999999999
$ws.Range("I89").Value = 1638.0526
$ws.Range("K89").Value = 8190.263
$ws.Range("M89").Value = -2574.263
$ws.Range("H94").Value = 1133.6333
$ws.Range("I94").Value = 612.9474
$ws.Range("K94").Value = 612.9474
$ws.Range("M94").Value = -161.9474
$ws.Range("H113").Value = 4874.875
$ws.Range("I113").Value = 4874.875
$ws.Range("K113").Value = 4874.875
$ws.Range("M113").Value = -2704.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3965.7334
$ws.Range("J31").Value = 5537.92
$ws.Range("L31").Value = 5537.92
$ws.Range("N31").Value = -6127.92
$ws.Range("H34").Value = 3965.7334
$ws.Range("J34").Value = 5537.92
$ws.Range("L34").Value = 5537.92
$ws.Range("N34").Value = -5941.92

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 5000
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15228
$ws.Range("H55").Value = 1288628.2
$ws.Range("J55").Value = 4999.5
$ws.Range("L55").Value = 14998.5
$ws.Range("N55").Value = -15352.5
$ws.Range("H56").Value = 22346.6
$ws.Range("I56").Value = 22346.6
$ws.Range("K56").Value = 22346.6
$ws.Range("M56").Value = -21816.6
$ws.Range("H131").Value = 1289.5186
$ws.Range("I131").Value = 936.625
$ws.Range("J131").Value = 1438.1052
$ws.Range("K131").Value = 2809.875
$ws.Range("L131").Value = 4314.3156
$ws.Range("M131").Value = 2230.125
$ws.Range("N131").Value = -14394.3156
$ws.Range("H134").Value = 1849.125
$ws.Range("J134").Value = 4033
$ws.Range("L134").Value = 12099
$ws.Range("N134").Value = -22239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 16179.8
$ws.Range("I10").Value = 20099.75
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 20099.75
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = -19930.75
$ws.Range("N10").Value = -838
$ws.Range("H14").Value = 7509250
$ws.Range("I14").Value = 7509250
$ws.Range("K14").Value = 7509250
$ws.Range("M14").Value = -7509082
$ws.Range("H18").Value = 2512000
$ws.Range("I18").Value = 2512000
$ws.Range("K18").Value = 2512000
$ws.Range("M18").Value = -2511707
$ws.Range("H47").Value = 43499
$ws.Range("J47").Value = 43499
$ws.Range("L47").Value = 43499
$ws.Range("N47").Value = -44635
$ws.Range("H70").Value = 5964.857
$ws.Range("I70").Value = 5465.8887
$ws.Range("J70").Value = 6863
$ws.Range("K70").Value = 5465.8887
$ws.Range("L70").Value = 6863
$ws.Range("M70").Value = -5195.8887
$ws.Range("N70").Value = -7403
$ws.Range("H73").Value = 5964.857
$ws.Range("I73").Value = 5465.8887
$ws.Range("J73").Value = 6863
$ws.Range("K73").Value = 5465.8887
$ws.Range("L73").Value = 6863
$ws.Range("M73").Value = -4529.8887
$ws.Range("N73").Value = -8735
$ws.Range("H80").Value = 3919.8
$ws.Range("I80").Value = 3200
$ws.Range("K80").Value = 3200
$ws.Range("M80").Value = -2202
$ws.Range("H83").Value = 3919.8
$ws.Range("I83").Value = 3200
$ws.Range("K83").Value = 16000
$ws.Range("M83").Value = -11008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4703.625
$ws.Range("I7").Value = 3806.5625
$ws.Range("K7").Value = 3806.5625
$ws.Range("M7").Value = -3694.5625
$ws.Range("H16").Value = 1439.9667
$ws.Range("J16").Value = 2200
$ws.Range("L16").Value = 2200
$ws.Range("N16").Value = -2540
$ws.Range("H40").Value = 3688.6667
$ws.Range("I40").Value = 3072.8
$ws.Range("K40").Value = 3072.8
$ws.Range("M40").Value = -2936.8
$ws.Range("H42").Value = 1271748.8
$ws.Range("I42").Value = 5012000
$ws.Range("J42").Value = 24998.334
$ws.Range("K42").Value = 5012000
$ws.Range("L42").Value = 24998.334
$ws.Range("M42").Value = -5011437
$ws.Range("N42").Value = -26124.334
$ws.Range("H43").Value = 22994.445
$ws.Range("J43").Value = 22992.857
$ws.Range("L43").Value = 22992.857
$ws.Range("N43").Value = -23378.857
$ws.Range("H49").Value = 1271748.8
$ws.Range("I49").Value = 5012000
$ws.Range("J49").Value = 24998.334
$ws.Range("K49").Value = 5012000
$ws.Range("L49").Value = 24998.334
$ws.Range("M49").Value = -5011853
$ws.Range("N49").Value = -25292.334
$ws.Range("H61").Value = 4265.1
$ws.Range("I61").Value = 3850.111
$ws.Range("K61").Value = 3850.111
$ws.Range("M61").Value = -3648.111
$ws.Range("H82").Value = 4060.2666
$ws.Range("I82").Value = 4298
$ws.Range("K82").Value = 4298
$ws.Range("M82").Value = -3937
$ws.Range("H85").Value = 4060.2666
$ws.Range("I85").Value = 4298
$ws.Range("K85").Value = 4298
$ws.Range("M85").Value = -3050
$ws.Range("H93").Value = 2385.6667
$ws.Range("I93").Value = 2375.2727
$ws.Range("K93").Value = 2375.2727
$ws.Range("M93").Value = -1127.2727
$ws.Range("H100").Value = 3587.682
$ws.Range("I100").Value = 3129.4443
$ws.Range("J100").Value = 5649.75
$ws.Range("K100").Value = 3129.4443
$ws.Range("L100").Value = 5649.75
$ws.Range("M100").Value = -2588.4443
$ws.Range("N100").Value = -6731.75
$ws.Range("H113").Value = 4265.1
$ws.Range("I113").Value = 3850.111
$ws.Range("K113").Value = 3850.111
$ws.Range("M113").Value = -1680.111
$ws.Range("H122").Value = 4115.1333
$ws.Range("I122").Value = 3811.6365
$ws.Range("J122").Value = 4949.75
$ws.Range("K122").Value = 11434.9095
$ws.Range("L122").Value = 14849.25
$ws.Range("M122").Value = -8984.9095
$ws.Range("N122").Value = -19749.25
$ws.Range("H126").Value = 4703.625
$ws.Range("I126").Value = 3806.5625
$ws.Range("K126").Value = 11419.6875
$ws.Range("M126").Value = -8949.6875
$ws.Range("H136").Value = 3860.95
$ws.Range("I136").Value = 3967.7778
$ws.Range("J136").Value = 2899.5
$ws.Range("K136").Value = 11903.3334
$ws.Range("L136").Value = 8698.5
$ws.Range("M136").Value = -9353.3334
$ws.Range("N136").Value = -13798.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8342.857
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 8342.857
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H96").Value = 4696.3687
$ws.Range("I96").Value = 3994.2222
$ws.Range("K96").Value = 3994.2222
$ws.Range("M96").Value = -2621.2222
$ws.Range("H100").Value = 31563898
$ws.Range("I100").Value = 42084896
$ws.Range("J100").Value = 900.625
$ws.Range("K100").Value = 84169792
$ws.Range("L100").Value = 1801.25
$ws.Range("M100").Value = -84169251
$ws.Range("N100").Value = -2883.25
$ws.Range("H122").Value = 4883.2856
$ws.Range("I122").Value = 4883.2856
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14649.8568
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -12199.8568
$ws.Range("H126").Value = 5874.8096
$ws.Range("I126").Value = 5601.2144
$ws.Range("K126").Value = 16803.6432
$ws.Range("M126").Value = -14333.6432
$ws.Range("H136").Value = 5341.592
$ws.Range("I136").Value = 3176.5483
$ws.Range("K136").Value = 9529.644899999999
$ws.Range("M136").Value = -6979.644899999999
